# Generate Report for Handback
#
# The f0730ab5-7411-454b-917c-92d52dd70f04.md file has now been handed back
# (its latest handback xliff is in sync with en-US), so:
#   - its Status flips from "Ready for handoff" to
#     "Handed back: in sync with en-US" on every sheet that tracks it
#     (Overview, zh-cn, de-de)
#   - the per-language "Latest Handback DateTime" gets stamped with the
#     generation time of the new handback report
#   - the stale "version mismatch" Error Detail message is cleared now that
#     the handback is current

$wb = $excel.ActiveWorkbook

# --- Overview sheet --------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet -------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"
$zhcn.Range("K3").Value = "2016-08-26 06:47:02"
$zhcn.Range("P3").Value = ""

# --- de-de sheet ---------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Handed back: in sync with en-US"
$dede.Range("K3").Value = "2016-08-26 06:47:13"
$dede.Range("P3").Value = ""
